$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.727.44'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '2.780.72'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '356.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.99'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.556'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.48%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.26%  '
$ws.Range("E11").Value = '  +3.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0844'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.51'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("D15").Value = '3.219.34'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '2.796.30'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '51.685.50'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("E20").Value = '  -3.51%  '
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.18'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.35'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.74'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.35'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.166'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +14.90%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.13'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.21%  '
$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.94'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.83'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0447'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -8.46%  '
$ws.Range("E35").Value = '  -2.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.18'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -7.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.85'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.56%  '
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("E40").Value = '  -4.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.55'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.81'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.68'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.91%  '
$ws.Range("D46").Value = '2.090.08'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("E47").Value = '  -2.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.952'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.09%  '
$ws.Range("E50").Value = '  -6.82%  '
$ws.Range("E51").Value = '  -7.18%  '
